$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Handback has happened and is in sync with en-US: update the Status cells
# everywhere they appear (Overview's per-locale columns, and each locale's
# own "Status" column).
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime refreshed for both locales.
$wsZhCn.Range("K2").Value = "2016-10-24 06:36:43"
$wsDeDe.Range("K2").Value = "2016-10-24 06:36:59"

# Error Detail cleared now that the handback is in sync.
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# Refresh the column widths that depend on the new cell contents.
$wsOverview.Range("E1:F2").EntireColumn.AutoFit()
$wsZhCn.Range("C1:C2").EntireColumn.AutoFit()
$wsZhCn.Range("P1:P2").EntireColumn.AutoFit()
$wsDeDe.Range("C1:C2").EntireColumn.AutoFit()
$wsDeDe.Range("P1:P2").EntireColumn.AutoFit()
